$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as plain text (e.g. "25.974.42", "213.95").
# Forcing the data range to Text format before writing keeps Excel from
# auto-converting numeric-looking strings (like "213.95") into real numbers,
# then ClearFormats restores the cells to their original (unstyled) look.
$ws.Range("D2:D51").NumberFormat = "@"

# Price (D column) updates
$ws.Range("D2").Value = "25.974.42"
$ws.Range("D3").Value = "1.633.16"
$ws.Range("D5").Value = "213.95"
$ws.Range("D8").Value = "0.252"
$ws.Range("D9").Value = "0.0624"
$ws.Range("D10").Value = "18.52"
$ws.Range("D12").Value = "1.859.61"
$ws.Range("D13").Value = "1.655.64"
$ws.Range("D15").Value = "0.531"
$ws.Range("D16").Value = "0.0₃0746"
$ws.Range("D17").Value = "25.983.90"
$ws.Range("D18").Value = "61.72"
$ws.Range("D20").Value = "190.38"
$ws.Range("D21").Value = "4.24"
$ws.Range("D22").Value = "9.57"
$ws.Range("D23").Value = "6.12"
$ws.Range("D24").Value = "0.132"
$ws.Range("D25").Value = "143.41"
$ws.Range("D28").Value = "6.78"
$ws.Range("D29").Value = "15.21"
$ws.Range("D31").Value = "0.0484"
$ws.Range("D32").Value = "3.16"
$ws.Range("D33").Value = "3.15"
$ws.Range("D36").Value = "1.132.86"
$ws.Range("D37").Value = "0.866"
$ws.Range("D41").Value = "98.57"
$ws.Range("D44").Value = "1.769.91"
$ws.Range("D46").Value = "55.10"
$ws.Range("D50").Value = "7.54"

$ws.Range("D2:D51").ClearFormats()

# Volume(1h) (E column) updates
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("E10").Value = "  -5.67%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("E15").Value = "  -2.62%  "
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("E33").Value = "  -4.21%  "
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("E35").Value = "  -2.47%  "
$ws.Range("E37").Value = "  -4.18%  "
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("E43").Value = "  -4.95%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("E51").Value = "  +0.14%  "
